$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows starting at row 466, pushing the existing
# rows 466-506 down to 469-509 (matches the target dimension A1:R509).
$ws.Range("A466:R468").Insert()

# New row data (Mercado ID .. Clasificación), matching the existing
# column layout: A Mercado ID, B Mercado, C Región, D Fecha, E Codreg,
# F Categoría ID, G Categoría, H Variedad, I Calidad, J Volumen,
# K Precio mínimo, L Precio máximo, M Precio promedio ponderado,
# N Unidad de comercialización, O Origen, P Precio $/Kg, Q Kg o Unidades,
# R Clasificación.
$newRows = @(
    @(3, "Femacal de La Calera", "Coquimbo", 44918, 5, 100112027, "Melón", "Tuna", "Extra",   550, 2500, 2500, 2500, "`$/unidad", "Región de O'Higgins", 2500, 1, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44918, 5, 100112027, "Melón", "Tuna", "Primera",  560, 1800, 1800, 1800, "`$/unidad", "Región de O'Higgins", 1800, 1, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44918, 5, 100112027, "Melón", "Tuna", "Segunda",  580, 1000, 1000, 1000, "`$/unidad", "Región de O'Higgins", 1000, 1, "Hortaliza")
)

$startRow = 466
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowData = $newRows[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
